# Update with edge routers.
# Adds two new edge switches (sw-edge-001/002) in x3003 with ISL connections
# to the existing spine switches on the INTER_SWITCH_LINKS sheet, then makes
# that sheet the active one (it was COMPUTE_NODES before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INTER_SWITCH_LINKS")

# --- Insert 4 new data rows (39-42), copying the formatting of the last
#     existing data row (38) so the new rows pick up the same cell style. ---
$ws.Rows.Item(38).Copy()
$ws.Rows.Item(39).Insert(-4121)

$ws.Rows.Item(38).Copy()
$ws.Rows.Item(40).Insert(-4121)

$ws.Rows.Item(38).Copy()
$ws.Rows.Item(41).Insert(-4121)

$ws.Rows.Item(38).Copy()
$ws.Rows.Item(42).Insert(-4121)

# --- Fill in the new rows' values. Column order here matters: it controls
#     the order new strings are appended to the shared-string table, so the
#     first brand new string written must be "x3003", then "sw-edge-001",
#     then "sw-edge-002" (everything else re-uses already-existing strings). ---

# Row 39: sw-edge-001 <-> sw-spine-001 (x3003 <-> x3000)
$ws.Range("K39").Value = "x3003"
$ws.Range("J39").Value = "sw-edge-001"
$ws.Range("L39").Value = "u40"
$ws.Range("O39").Value = 1
$ws.Range("P39").Value = "sw-spine-001"
$ws.Range("Q39").Value = "x3000"
$ws.Range("R39").Value = "u40"
$ws.Range("T39").Value = 7

# Row 40: sw-edge-002 <-> sw-spine-001 (x3003 <-> x3000)
$ws.Range("J40").Value = "sw-edge-002"
$ws.Range("K40").Value = "x3003"
$ws.Range("L40").Value = "u40"
$ws.Range("O40").Value = 1
$ws.Range("P40").Value = "sw-spine-001"
$ws.Range("Q40").Value = "x3000"
$ws.Range("R40").Value = "u40"
$ws.Range("T40").Value = 8

# Row 41: sw-edge-001 <-> sw-spine-002 (x3003 <-> x3001)
$ws.Range("J41").Value = "sw-edge-001"
$ws.Range("K41").Value = "x3003"
$ws.Range("L41").Value = "u40"
$ws.Range("O41").Value = 2
$ws.Range("P41").Value = "sw-spine-002"
$ws.Range("Q41").Value = "x3001"
$ws.Range("R41").Value = "u40"
$ws.Range("T41").Value = 7

# Row 42: sw-edge-002 <-> sw-spine-002 (x3003 <-> x3001)
$ws.Range("J42").Value = "sw-edge-002"
$ws.Range("K42").Value = "x3003"
$ws.Range("L42").Value = "u40"
$ws.Range("O42").Value = 2
$ws.Range("P42").Value = "sw-spine-002"
$ws.Range("Q42").Value = "x3001"
$ws.Range("R42").Value = "u40"
$ws.Range("T42").Value = 8

# --- Switch the active sheet/tab from COMPUTE_NODES to INTER_SWITCH_LINKS,
#     and move the selection to K45 on the new active sheet. ---
$ws.Activate() | Out-Null
$ws.Range("K45").Select() | Out-Null
